$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = -70351.28293002227
$ws.Range("B7").Value = 11132708.25685823
$ws.Range("B8").Value = 25408340.88108395
$ws.Range("B10").Value = 2756026.145368536

$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("L11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("K12").Value = 80.29914934735042
$ws.Range("L12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("P14").Value = 135.4597561231036
$ws.Range("Q14").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("J20").Value = 124.5190384721106
$ws.Range("P20").Value = 135.4597561231036
$ws.Range("Q20").Value = 150.3839754851235
$ws.Range("L21").Value = 0
$ws.Range("L22").Value = 90.4687457914608
$ws.Range("M22").Value = 92.09541281912071
$ws.Range("N22").Value = 81.96869489115805
$ws.Range("K23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("K27").Value = 80.29914934735042
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("P29").Value = 135.4597561231036
$ws.Range("Q29").Value = 150.3839754851235
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 96.22962838366004
$ws.Range("L32").Value = 0
$ws.Range("K33").Value = 80.29914934735042
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 65.92768427608706
$ws.Range("N34").Value = 0
$ws.Range("K36").Value = 80.29914934735042
$ws.Range("L37").Value = 90.4687457914608
$ws.Range("M37").Value = 92.09541281912071
$ws.Range("O37").Value = 96.22962838366004
$ws.Range("Q39").Value = 0
$ws.Range("J42").Value = 93.17061249236157
$ws.Range("K43").Value = 0
$ws.Range("O43").Value = 96.22962838366004
$ws.Range("Q44").Value = 150.3839754851235
$ws.Range("L46").Value = 90.4687457914608
$ws.Range("Q46").Value = 65.34295837775146

$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("L11").Value = 130.6648563030561
$ws.Range("O11").Value = 117.8828208804077
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 61.18167021676314
$ws.Range("O12").Value = 57.81213424001893
$ws.Range("L14").Value = 130.6648563030561
$ws.Range("N14").Value = 110.5750244233121
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 150.3839754851235
$ws.Range("M15").Value = 51.84373129681028
$ws.Range("N15").Value = 38.66169381481656
$ws.Range("O15").Value = 57.81213424001893
$ws.Range("Q15").Value = 94.49434172313325
$ws.Range("L16").Value = 90.4687457914608
$ws.Range("M16").Value = 92.09541281912071
$ws.Range("N16").Value = 81.96869489115805
$ws.Range("O16").Value = 96.22962838366004
$ws.Range("M19").Value = 92.09541281912071
$ws.Range("J20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("L21").Value = 61.18167021676314
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("K23").Value = 135.370731907559
$ws.Range("Q23").Value = 150.3839754851235
$ws.Range("J24").Value = 93.17061249236157
$ws.Range("K24").Value = 80.29914934735042
$ws.Range("L24").Value = 61.18167021676314
$ws.Range("M24").Value = 51.84373129681028
$ws.Range("O24").Value = 57.81213424001893
$ws.Range("P24").Value = 65.92768427608706
$ws.Range("Q24").Value = 94.49434172313325
$ws.Range("L25").Value = 90.4687457914608
$ws.Range("M25").Value = 92.09541281912071
$ws.Range("P25").Value = 101.5955875616828
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 61.18167021676314
$ws.Range("M27").Value = 51.84373129681028
$ws.Range("N27").Value = 38.66169381481656
$ws.Range("O27").Value = 57.81213424001893
$ws.Range("P27").Value = 65.92768427608706
$ws.Range("M28").Value = 92.09541281912071
$ws.Range("N28").Value = 81.96869489115805
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("N31").Value = 81.96869489115805
$ws.Range("O31").Value = 0
$ws.Range("L32").Value = 130.6648563030561
$ws.Range("K33").Value = 0
$ws.Range("O33").Value = 57.81213424001893
$ws.Range("P33").Value = 0
$ws.Range("N34").Value = 81.96869489115805
$ws.Range("K36").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("Q39").Value = 94.49434172313325
$ws.Range("J42").Value = 0
$ws.Range("K43").Value = 94.30397654773019
$ws.Range("O43").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("Q46").Value = 61.14583096471014

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B5").Value = 281822.62736734
$ws.Range("B6").Value = 257407.2830265637
$ws.Range("B7").Value = 178032.8048269036
$ws.Range("B8").Value = 255219.2967190486
$ws.Range("B9").Value = 191440.6432236616
$ws.Range("B10").Value = 191590.6143084067
$ws.Range("B11").Value = 210608.4641312597
$ws.Range("B12").Value = 273667.3660632239
$ws.Range("B13").Value = 246349.5939647874
$ws.Range("B14").Value = 192943.1099506649
$ws.Range("B15").Value = 188982.3429356603
$ws.Range("B16").Value = 187398.2196802687

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("E2").Value = 39162.55191854898
$ws.Range("F2").Value = 36023.43621759204
$ws.Range("G2").Value = 25818.14616335002
$ws.Range("H2").Value = 35742.12369234009
$ws.Range("I2").Value = 27542.01110007605
$ws.Range("J2").Value = 27561.29309668612
$ws.Range("K2").Value = 30006.44521676723
$ws.Range("L2").Value = 38114.01832230548
$ws.Range("M2").Value = 34601.7333382208
$ws.Range("N2").Value = 27735.1853935479
$ws.Range("O2").Value = 27225.94392019016
$ws.Range("P2").Value = 27022.27093021123
$ws.Range("E3").Value = 133100.0000000001
$ws.Range("C4").Value = 48378.33248915088
$ws.Range("E4").Value = 16387.55325273026
$ws.Range("F4").Value = 13248.43755177332
$ws.Range("G4").Value = 3043.1474975313
$ws.Range("H4").Value = 12967.12502652138
$ws.Range("I4").Value = 4767.012434257326
$ws.Range("J4").Value = 4786.2944308674
$ws.Range("K4").Value = 7231.446550948505
$ws.Range("L4").Value = 15339.01965648677
$ws.Range("M4").Value = 11826.73467240208
$ws.Range("N4").Value = 4960.186727729181
$ws.Range("O4").Value = 4450.945254371439
$ws.Range("P4").Value = 4247.272264392515
$ws.Range("B6").Value = -33627.59999999999
$ws.Range("C6").Value = -33627.6
$ws.Range("E6").Value = -113936.6134391897
$ws.Range("H6").Value = 19163.38656081039
$ws.Range("I6").Value = 19163.3865608104
